$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; unprotect to allow edits, then re-protect at the end.
$ws.Unprotect("D382")

# Update the "as of" date in the confidential disclaimer text (cell A9).
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + "`n" + "Model holdings provided as of 2021-05-12 for illustrative purposes only and are subject to change."

# Update Weight (column D) and Percent Change (column E) values for rows 2-6.
$ws.Range("D2").Value = 0.2584154364642092
$ws.Range("E2").Value = -0.02517944628004987

$ws.Range("D3").Value = 0.2519571370232985
$ws.Range("E3").Value = -0.01024811218985977

$ws.Range("D4").Value = 0.245224390562171
$ws.Range("E4").Value = -0.0239676580999133

$ws.Range("D5").Value = 0.2444030359503214
$ws.Range("E5").Value = -0.01959247648902829

$ws.Range("E6").Value = -0.01975475769409418

# Restore sheet protection.
$ws.Protect("D382")
